$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 27 (the "cholesterol HDL" duplicate row), shifting rows below it up by one.
$ws.Rows.Item(27).Delete()

# Update the view's selection / scroll position to match the target state.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("E49").Select()
